# Insert a new row at row 186, pushing existing rows 186:218 down to 187:219,
# and populate the new row with the latest weekly price record for Perejil.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(186).Insert()

$ws.Range("A186").Value = 9
$ws.Range("B186").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C186").Value = "Metropolitana"
$ws.Range("D186").Value = 44474
$ws.Range("D186").NumberFormat = $ws.Range("D187").NumberFormat
$ws.Range("E186").Value = 13
$ws.Range("F186").Value = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 106
$ws.Range("K186").Value = 8000
$ws.Range("L186").Value = 10000
$ws.Range("M186").Value = 9000
$ws.Range("N186").Value = "`$/docena de atados"
$ws.Range("O186").Value = "Región Metropolitana"
$ws.Range("P186").Value = 3000
$ws.Range("Q186").Value = 3
$ws.Range("R186").Value = "Hortaliza"
